$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the schedule values for trial 1 (row 2)
$ws.Range("D2").Value = 6
$ws.Range("F2").Value = -3
$ws.Range("H2").Value = 46

# Move the active selection to H5, matching the saved cursor position
$ws.Range("H5").Select()
